# Daily refresh of the cryptos price/volume table (GitHub Actions bot).
# Column D = Price, Column E = Volume(1h) change.
# Some "D" values look numeric (e.g. "19.73"); Excel would otherwise silently
# coerce them from text to a real number on assignment, which changes the
# cell's stored type/formatting versus the source data (plain text cells).
# Prefix those with an apostrophe to force text entry, then reset the style
# to "Normal" so the forced text-format style doesn't stick to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.666.67"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "1.596.57"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'211.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "'0.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").Value = "'19.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").Value = "1.820.68"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").Value = "1.597.76"
$ws.Range("E13").Value = "  -1.71%  "
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "'65.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "26.640.74"
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "'209.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D23").Value = "'2.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").Value = "'8.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'146.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D27").Value = "'7.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.81%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("E32").Value = "  -3.23%  "
$ws.Range("D33").Value = "'0.670"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.56%  "
$ws.Range("D34").Value = "'2.92"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "1.293.40"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("D36").Value = "'2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").Value = "'0.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'5.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").Value = "'0.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "'2.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "'63.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.49%  "
$ws.Range("D45").Value = "1.733.01"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").Value = "'89.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'1.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").Value = "'0.860"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").Value = "'0.0986"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.21%  "
$ws.Range("D50").Value = "'0.0504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").Value = "'7.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.74%  "
